$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the column headers in row 1: "_old" -> "_FV2410", "_new" -> "_FV2504" ---
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value() -replace "_old$", "_FV2410")
}

foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value() -replace "_new$", "_FV2504")
}

# --- 2) Turn the data range into an Excel Table (ListObject) ---
$rng = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# --- 3) Freeze the header row (pane split after row 1) ---
$ws.Range("A2").Select()
$excel.ActiveWindow().FreezePanes = $true
